$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("estimates")

# Insert a new blank row at 116, pushing the existing rows 116-118 down to 117-119.
$ws.Rows("116:116").Insert()

# Fill in the new row's data. Column order below matches the order the
# strings were actually entered (so new shared-string indices line up):
# G, I, K, C first (these introduce brand-new shared strings), then the rest.
$ws.Range("G116").Value = "emory.aggregated"
$ws.Range("I116").Value = "section.3.processing"
$ws.Range("K116").Value = "Used pull function to aggregate data from county to MSA and state. Note that we are only aggregating data that we have so there may be counties missing from a state (in Alaska for example)."
$ws.Range("K116").HorizontalAlignment = -4108
$ws.Range("C116").Value = "MSA; State"

$ws.Range("A116").Value = "Section3"
$ws.Range("B116").Value = "proportion.msm"
$ws.Range("D116").Value = "yes"
$ws.Range("E116").Value = 2013
$ws.Range("F116").Value = "sex*"
$ws.Range("H116").Value = "emory"
$ws.Range("J116").Value = "proportion"

# No Link to Data for this new row - remove the blank cell that the
# row-insert auto-created (with inherited formatting) in column L.
$ws.Range("L116").Clear()

# Restore the auto filter / filter database range to cover the extra row
# (Excel grows these by one row when a row is inserted above their bottom edge).
$ws.AutoFilterMode = $false
$ws.Range("A1:L120").AutoFilter()
foreach ($n in $wb.Names) {
  if ($n.Name -eq "estimates!_FilterDatabase") {
    $n.RefersTo = "=estimates!`$A`$1:`$L`$120"
  }
}

# Reflect the saved selection/active cell position from the source edit.
$ws.Range("D132").Select()

Write-Host "Row inserted and populated"
